$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '52.693.54'
$ws.Range("E2").Value = '  -12.80%  '

$ws.Range("D3").Value = '2.311.92'
$ws.Range("E3").Value = '  -20.12%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").Value = '''441.48'
$ws.Range("E5").Value = '  -15.94%  '

$ws.Range("D6").Value = '''122.21'
$ws.Range("E6").Value = '  -13.49%  '

$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").Value = '''0.467'
$ws.Range("E8").Value = '  -14.74%  '

$ws.Range("D9").Value = '2.317.09'
$ws.Range("E9").Value = '  -20.05%  '

$ws.Range("D10").Value = '''5.30'
$ws.Range("E10").Value = '  -11.29%  '

$ws.Range("D11").Value = '''0.0877'
$ws.Range("E11").Value = '  -17.92%  '

$ws.Range("D12").Value = '''0.302'
$ws.Range("E12").Value = '  -15.49%  '

$ws.Range("E13").Value = '  -5.96%  '

$ws.Range("D14").Value = '53.022.40'
$ws.Range("E14").Value = '  -12.30%  '

$ws.Range("D15").Value = '''18.86'
$ws.Range("E15").Value = '  -16.76%  '

$ws.Range("E16").Value = '  -15.26%  '

$ws.Range("D17").Value = '2.333.62'
$ws.Range("E17").Value = '  -19.58%  '

$ws.Range("D18").Value = '''3.94'
$ws.Range("E18").Value = '  -20.78%  '

$ws.Range("D19").Value = '''297.84'
$ws.Range("E19").Value = '  -15.46%  '

$ws.Range("D20").Value = '''8.91'
$ws.Range("E20").Value = '  -23.28%  '

$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").Value = '''5.61'
$ws.Range("E22").Value = '  -1.59%  '

$ws.Range("D23").Value = '''5.13'
$ws.Range("E23").Value = '  -21.78%  '

$ws.Range("D24").Value = '''53.68'
$ws.Range("E24").Value = '  -16.67%  '

$ws.Range("D25").Value = '''0.148'
$ws.Range("E25").Value = '  -16.87%  '

$ws.Range("D26").Value = '''0.363'
$ws.Range("E26").Value = '  -19.71%  '

$ws.Range("B27").Value = 'USDe'
$ws.Range("C27").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D27").Value = '''0.997'
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '''6.81'
$ws.Range("E28").Value = '  -12.72%  '

$ws.Range("D29").Value = '0.0₃0673'
$ws.Range("E29").Value = '  -19.31%  '

$ws.Range("D30").Value = '''141.25'
$ws.Range("E30").Value = '  -6.01%  '

$ws.Range("D31").Value = '''16.90'
$ws.Range("E31").Value = '  -13.61%  '

$ws.Range("E32").Value = '  -19.96%  '

$ws.Range("D33").Value = '''4.75'
$ws.Range("E33").Value = '  -14.66%  '

$ws.Range("D34").Value = '''0.828'
$ws.Range("E34").Value = '  -16.74%  '

$ws.Range("D35").Value = '''3.42'
$ws.Range("E35").Value = '  -20.64%  '

$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = '''0.994'
$ws.Range("E36").Value = '  -0.31%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.993'
$ws.Range("E37").Value = '  -17.16%  '

$ws.Range("D38").Value = '''31.72'
$ws.Range("E38").Value = '  -15.73%  '

$ws.Range("E39").Value = '  -1.65%  '

$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").Value = '''0.557'
$ws.Range("E40").Value = '  -13.66%  '

$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '''0.0503'
$ws.Range("E41").Value = '  -13.12%  '

$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '''3.13'
$ws.Range("E42").Value = '  -15.57%  '

$ws.Range("D43").Value = '1.908.69'
$ws.Range("E43").Value = '  -16.39%  '

$ws.Range("D44").Value = '''1.19'
$ws.Range("E44").Value = '  -18.82%  '

$ws.Range("D45").Value = '''4.26'
$ws.Range("E45").Value = '  -13.48%  '

$ws.Range("D46").Value = '''0.0822'
$ws.Range("E46").Value = '  -10.32%  '

$ws.Range("D47").Value = '''0.0205'
$ws.Range("E47").Value = '  -13.25%  '

$ws.Range("D48").Value = '''15.57'
$ws.Range("E48").Value = '  -23.16%  '

$ws.Range("D49").Value = '''4.58'
$ws.Range("E49").Value = '  -5.16%  '

$ws.Range("D50").Value = '''4.44'
$ws.Range("E50").Value = '  -13.33%  '

$ws.Range("D51").Value = '''15.01'
$ws.Range("E51").Value = '  -17.36%  '
